$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the Game 1..6 / Good Session headers (keep Notes_t as-is)
$ws.Range("G1").Value = "Game.1"
$ws.Range("H1").Value = "Game.2"
$ws.Range("I1").Value = "Game.3"
$ws.Range("J1").Value = "Game.4"
$ws.Range("K1").Value = "Game.5"
$ws.Range("L1").Value = "Game.6"
$ws.Range("M1").Value = "goodSession"

# Move the active selection/view to G2 (frozen pane scrolled back to the top)
$ws.Range("G2").Select()
